$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 4 updates
$ws.Range("B4").Value = 0.75
$ws.Range("G4").Value = 35
$ws.Range("H4").Value = 35
$ws.Range("J4").Value = 50
$ws.Range("K4").Value = 50

# Row 6 updates
$ws.Range("B6").Value = 0.75
$ws.Range("G6").Value = 20
$ws.Range("H6").Value = 20

# Update the active selection to B4 (matches recorded sheetView selection)
$ws.Activate()
$ws.Range("B4").Select()
